$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sector_Fuels")
$ws.Activate()
